$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("room_type", "price", "sum", 0.8439496330405888),
    @("room_type", "beds", "sum", 0.6924646033478833),
    @("room_type", "number_of_reviews", "sum", 0.6403318649977773),
    @("room_type", "price", "max", 0.5234970057746275),
    @("neighbourhood", "number_of_reviews", "avg", 0.3784619634968182),
    @("room_type", "number_of_reviews", "max", 0.3764841538949056),
    @("room_type", "number_of_reviews", "avg", 0.3459132348469507),
    @("neighbourhood", "number_of_reviews", "max", 0.2801054798777837),
    @("neighbourhood", "price", "max", 0.2597512415111768),
    @("neighbourhood", "price", "sum", 0.2572740350232723),
    @("neighbourhood", "beds", "max", 0.23294698605459),
    @("zipcode", "number_of_reviews", "avg", 0.2266450825819268),
    @("neighbourhood", "beds", "sum", 0.2239485353969972),
    @("room_type", "price", "avg", 0.197771542884957),
    @("zipcode", "number_of_reviews", "max", 0.173818684796122),
    @("zipcode", "number_of_reviews", "sum", 0.1709315449053815),
    @("zipcode", "price", "sum", 0.1704145754511847),
    @("zipcode", "price", "max", 0.1672195758327164),
    @("zipcode", "beds", "sum", 0.146114404950847),
    @("zipcode", "price", "avg", 0.1456157995716201),
    @("zipcode", "beds", "avg", 0.1433027386152141),
    @("zipcode", "beds", "max", 0.138432003892941),
    @("room_type", "beds", "max", 0.1263209819360183),
    @("neighbourhood", "beds", "avg", 0.08299289533362203),
    @("room_type", "beds", "avg", 0.07981191646023601),
    @("neighbourhood", "price", "avg", 0.07358113058437878),
    @("neighbourhood", "number_of_reviews", "sum", 0.07110365005260268)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $data[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $data[$i][2]
    $ws.Cells.Item($r, 4).Value2 = $data[$i][3]
}

$ws.Application.GoTo($ws.Range("J21"), $true)